# Weekly update: insert a new price record for Rabanito (row 147) and
# push the existing rows 147-181 down by one (they keep their original
# values, now living one row lower).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 147 - this shifts rows 147:181 down to
# 148:182 and keeps all their existing data/formatting intact.
$ws.Rows.Item(147).Insert()

# Populate the newly inserted row 147 with the new record.
$ws.Cells.Item(147, 1).Value = 9
$ws.Cells.Item(147, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(147, 3).Value = "Metropolitana"
$ws.Cells.Item(147, 4).Value = 44511
$ws.Cells.Item(147, 5).Value = 13
$ws.Cells.Item(147, 6).Value = 300000001
$ws.Cells.Item(147, 7).Value = "Rabanito"
$ws.Cells.Item(147, 8).Value = "Sin especificar"
$ws.Cells.Item(147, 9).Value = "Primera"
$ws.Cells.Item(147, 10).Value = 7900
$ws.Cells.Item(147, 11).Value = 2500
$ws.Cells.Item(147, 12).Value = 3000
$ws.Cells.Item(147, 13).Value = 2747
$ws.Cells.Item(147, 14).Value = "$/cien unidades (volumen en unidades)"
$ws.Cells.Item(147, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(147, 16).Value = 27
$ws.Cells.Item(147, 17).Value = 100
$ws.Cells.Item(147, 18).Value = "Hortaliza"

"done"
